$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row (row 1) from Russian labels to the new English identifiers
$ws.Range("A1").Value = "№"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "temp_in"
$ws.Range("D1").Value = "temp_out"
$ws.Range("E1").Value = "press_in"
$ws.Range("F1").Value = "press_out"
$ws.Range("G1").Value = "volume"
$ws.Range("H1").Value = "steps"
$ws.Range("I1").Value = "molar_flow_in"

# Update the selection shown in the sheet view
$ws.Range("B1:I1").Select()
